$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 341, shifting existing rows 341:361 down to 343:363
$ws.Rows("341:342").Insert()

# --- Row 341 (new): Femacal de La Calera / Coquimbo / Zapallo / Camote / 1a nueva(o) ---
$ws.Range("A341").Value = 3
$ws.Range("B341").Value = "Femacal de La Calera"
$ws.Range("C341").Value = "Coquimbo"
$ws.Range("D341").Value = 44516
$ws.Range("E341").Value = 5
$ws.Range("F341").Value = 100112045
$ws.Range("G341").Value = "Zapallo"
$ws.Range("H341").Value = "Camote"
$ws.Range("I341").Value = "1a nueva(o)"
$ws.Range("J341").Value = 160
$ws.Range("K341").Value = 600
$ws.Range("L341").Value = 600
$ws.Range("M341").Value = 600
$ws.Range("N341").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O341").Value = "Provincia de Talca"
$ws.Range("P341").Value = 600
$ws.Range("Q341").Value = 1
$ws.Range("R341").Value = "Hortaliza"

# --- Row 342 (new): Femacal de La Calera / Coquimbo / Zapallo / Paine / 1a nueva(o) ---
$ws.Range("A342").Value = 3
$ws.Range("B342").Value = "Femacal de La Calera"
$ws.Range("C342").Value = "Coquimbo"
$ws.Range("D342").Value = 44516
$ws.Range("E342").Value = 5
$ws.Range("F342").Value = 100112045
$ws.Range("G342").Value = "Zapallo"
$ws.Range("H342").Value = "Paine"
$ws.Range("I342").Value = "1a nueva(o)"
$ws.Range("J342").Value = 180
$ws.Range("K342").Value = 200
$ws.Range("L342").Value = 200
$ws.Range("M342").Value = 200
$ws.Range("N342").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O342").Value = "Provincia de Talca"
$ws.Range("P342").Value = 200
$ws.Range("Q342").Value = 1
$ws.Range("R342").Value = "Hortaliza"

Write-Host "Done applying edits"
